# A new weekly price record was added to the "Femacal de La Calera -
# Zanahoria" series. In the source system the rows are ordered (most
# recent first), so the new record lands at row 534 and every existing
# record from the old row 534 onward shifts down by one row.
#
# Insert a new row at position 534 (this pushes rows 534..636 down to
# 535..637 and keeps their values/formatting intact).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(534).Insert()

# Populate the newly inserted row 534 with the new record's data.
$ws.Cells.Item(534, 1).Value  = 3
$ws.Cells.Item(534, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(534, 3).Value  = "Coquimbo"
$ws.Cells.Item(534, 4).Value  = 45209
$ws.Cells.Item(534, 5).Value  = 5
$ws.Cells.Item(534, 6).Value  = 100114013
$ws.Cells.Item(534, 7).Value  = "Zanahoria"
$ws.Cells.Item(534, 8).Value  = "Sin especificar"
$ws.Cells.Item(534, 9).Value  = "Primera"
$ws.Cells.Item(534, 10).Value = 350
$ws.Cells.Item(534, 11).Value = 6500
$ws.Cells.Item(534, 12).Value = 7000
$ws.Cells.Item(534, 13).Value = 6771
$ws.Cells.Item(534, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(534, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(534, 16).Value = 339
$ws.Cells.Item(534, 17).Value = 20
$ws.Cells.Item(534, 18).Value = "Hortaliza"
